# Daily update at 8 AM UTC
# Adds a new row to the "Wins Over Time" tracking sheet with the next
# day's date plus the latest Chase/Bryce/Zach counts. The newest row's
# date cell is formatted as "YYYY-MM-DD" (no time) to highlight it as the
# latest entry; once a new row supersedes it, the old row reverts to the
# regular "YYYY-MM-DD HH:MM:SS" date format used by every other row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the current last used row in column A (the date column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
$newRow = $lastRow + 1

# The previous last row reverts to the regular date format used by every
# other historical row (instead of the "latest row" highlight format).
$ws.Range("A" + $lastRow).NumberFormat = $ws.Range("A" + ($lastRow - 1)).NumberFormat

# New day's data.
$ws.Cells.Item($newRow, 1).Value = 45836
$ws.Cells.Item($newRow, 2).Value = 414
$ws.Cells.Item($newRow, 3).Value = 408
$ws.Cells.Item($newRow, 4).Value = 421

# The newly appended row takes on the "latest row" highlight date format
# that the old last row previously had.
$ws.Range("A" + $newRow).NumberFormat = "YYYY-MM-DD"
